$d = $word.ActiveDocument

# 1. Delete the third paragraph entirely (the one holding the Doraemon picture/drawing).
$picPara = $d.Paragraphs.Item(3)
$picPara.Range.Delete()

# 2. Replace the second paragraph (currently styled "Title" with text "Doraemon ")
#    with a plain paragraph (no pStyle) containing the new sentence.
$titlePara = $d.Paragraphs.Item(2)
$titlePara.Range.InsertXML("<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r><w:t xml:space='preserve'>Doraemon is a cartoon character. </w:t></w:r></w:p>") | Out-Null
